$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 156.656447
$ws.Range("H2").Value = 469.969341
$ws.Range("I2").Value = 0.0671576211124673
$ws.Range("J2").Value = 0.0671576211124673
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.398836333333334
$ws.Range("N2").Value = 7.196509000000001
$ws.Range("O2").Value = 0.09386760623633866
$ws.Range("P2").Value = 0.09386760623633865
$ws.Range("Q2").Value = 375.7931769145077
$ws.Range("R2").Value = 3382.138592230569
$ws.Range("S2").Value = 0.006303925134354305
$ws.Range("T2").Value = 0.006303925134354304

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 156.656447
$ws.Range("H3").Value = 469.969341
$ws.Range("I3").Value = 0.0671576211124673
$ws.Range("J3").Value = 0.0671576211124673
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.682092999999999
$ws.Range("N3").Value = 17.046279
$ws.Range("O3").Value = 0.2223430006085962
$ws.Range("P3").Value = 0.2223430006085962
$ws.Range("Q3").Value = 890.1365009035708
$ws.Range("R3").Value = 8011.228508132138
$ws.Range("S3").Value = 0.01493202699188119
$ws.Range("T3").Value = 0.01493202699188119

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 156.656447
$ws.Range("H4").Value = 469.969341
$ws.Range("I4").Value = 0.0671576211124673
$ws.Range("J4").Value = 0.0671576211124673
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.17451333333334
$ws.Range("N4").Value = 36.52354
$ws.Range("O4").Value = 0.4763944950360188
$ws.Range("P4").Value = 0.4763944950360188
$ws.Range("Q4").Value = 1907.216002754127
$ws.Range("R4").Value = 17164.94402478714
$ws.Range("S4").Value = 0.03199352099769413
$ws.Range("T4").Value = 0.03199352099769413

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 156.656447
$ws.Range("H5").Value = 469.969341
$ws.Range("I5").Value = 0.0671576211124673
$ws.Range("J5").Value = 0.0671576211124673
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.300086333333334
$ws.Range("N5").Value = 15.900259
$ws.Range("O5").Value = 0.2073948981190463
$ws.Range("P5").Value = 0.2073948981190463
$ws.Range("Q5").Value = 830.2926937732576
$ws.Range("R5").Value = 7472.634243959319
$ws.Range("S5").Value = 0.01392814798853767
$ws.Range("T5").Value = 0.01392814798853767

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 42.300692
$ws.Range("H6").Value = 126.902076
$ws.Range("I6").Value = 0.01813403725498241
$ws.Range("J6").Value = 0.01813403725498241
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.398836333333334
$ws.Range("N6").Value = 7.196509000000001
$ws.Range("O6").Value = 0.09386760623633866
$ws.Range("P6").Value = 0.09386760623633865
$ws.Range("Q6").Value = 101.4724368947427
$ws.Range("R6").Value = 913.2519320526841
$ws.Range("S6").Value = 0.001702198668525784
$ws.Range("T6").Value = 0.001702198668525784

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 42.300692
$ws.Range("H7").Value = 126.902076
$ws.Range("I7").Value = 0.01813403725498241
$ws.Range("J7").Value = 0.01813403725498241
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.682092999999999
$ws.Range("N7").Value = 17.046279
$ws.Range("O7").Value = 0.2223430006085962
$ws.Range("P7").Value = 0.2223430006085962
$ws.Range("Q7").Value = 240.3564659083559
$ws.Range("R7").Value = 2163.208193175204
$ws.Range("S7").Value = 0.00403197625642086
$ws.Range("T7").Value = 0.00403197625642086

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 42.300692
$ws.Range("H8").Value = 126.902076
$ws.Range("I8").Value = 0.01813403725498241
$ws.Range("J8").Value = 0.01813403725498241
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 12.17451333333334
$ws.Range("N8").Value = 36.52354
$ws.Range("O8").Value = 0.4763944950360188
$ws.Range("P8").Value = 0.4763944950360188
$ws.Range("Q8").Value = 514.9903387632268
$ws.Range("R8").Value = 4634.913048869041
$ws.Range("S8").Value = 0.008638955521051697
$ws.Range("T8").Value = 0.008638955521051695

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 42.300692
$ws.Range("H9").Value = 126.902076
$ws.Range("I9").Value = 0.01813403725498241
$ws.Range("J9").Value = 0.01813403725498241
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.300086333333334
$ws.Range("N9").Value = 15.900259
$ws.Range("O9").Value = 0.2073948981190463
$ws.Range("P9").Value = 0.2073948981190463
$ws.Range("Q9").Value = 224.1973195597427
$ws.Range("R9").Value = 2017.775876037684
$ws.Range("S9").Value = 0.003760906808984066
$ws.Range("T9").Value = 0.003760906808984066

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2110.189616
$ws.Range("H10").Value = 6330.568848
$ws.Range("I10").Value = 0.9046248489651427
$ws.Range("J10").Value = 0.9046248489651426
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.398836333333334
$ws.Range("N10").Value = 7.196509000000001
$ws.Range("O10").Value = 0.09386760623633866
$ws.Range("P10").Value = 0.09386760623633865
$ws.Range("Q10").Value = 5061.999521083515
$ws.Range("R10").Value = 45557.99568975164
$ws.Range("S10").Value = 0.08491496911426735
$ws.Range("T10").Value = 0.08491496911426732

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2110.189616
$ws.Range("H11").Value = 6330.568848
$ws.Range("I11").Value = 0.9046248489651427
$ws.Range("J11").Value = 0.9046248489651426
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.682092999999999
$ws.Range("N11").Value = 17.046279
$ws.Range("O11").Value = 0.2223430006085962
$ws.Range("P11").Value = 0.2223430006085962
$ws.Range("Q11").Value = 11990.29364574629
$ws.Range("R11").Value = 107912.6428117166
$ws.Range("S11").Value = 0.201137003344008
$ws.Range("T11").Value = 0.201137003344008

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2110.189616
$ws.Range("H12").Value = 6330.568848
$ws.Range("I12").Value = 0.9046248489651427
$ws.Range("J12").Value = 0.9046248489651426
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 12.17451333333334
$ws.Range("N12").Value = 36.52354
$ws.Range("O12").Value = 0.4763944950360188
$ws.Range("P12").Value = 0.4763944950360188
$ws.Range("Q12").Value = 25690.53161585355
$ws.Range("R12").Value = 231214.7845426819
$ws.Range("S12").Value = 0.430958298119784
$ws.Range("T12").Value = 0.4309582981197839

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2110.189616
$ws.Range("H13").Value = 6330.568848
$ws.Range("I13").Value = 0.9046248489651427
$ws.Range("J13").Value = 0.9046248489651426
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 5.300086333333334
$ws.Range("N13").Value = 15.900259
$ws.Range("O13").Value = 0.2073948981190463
$ws.Range("P13").Value = 0.2073948981190463
$ws.Range("Q13").Value = 11184.18714450352
$ws.Range("R13").Value = 100657.6843005316
$ws.Range("S13").Value = 0.1876145783870835
$ws.Range("T13").Value = 0.1876145783870834

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 23.52144266666667
$ws.Range("H14").Value = 70.564328
$ws.Range("I14").Value = 0.01008349266740757
$ws.Range("J14").Value = 0.01008349266740757
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.398836333333334
$ws.Range("N14").Value = 7.196509000000001
$ws.Range("O14").Value = 0.09386760623633866
$ws.Range("P14").Value = 0.09386760623633865
$ws.Range("Q14").Value = 56.4240912812169
$ws.Range("R14").Value = 507.8168215309521
$ws.Range("S14").Value = 0.0009465133191912221
$ws.Range("T14").Value = 0.0009465133191912218

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 23.52144266666667
$ws.Range("H15").Value = 70.564328
$ws.Range("I15").Value = 0.01008349266740757
$ws.Range("J15").Value = 0.01008349266740757
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.682092999999999
$ws.Range("N15").Value = 17.046279
$ws.Range("O15").Value = 0.2223430006085962
$ws.Range("P15").Value = 0.2223430006085962
$ws.Range("Q15").Value = 133.651024726168
$ws.Range("R15").Value = 1202.859222535512
$ws.Range("S15").Value = 0.002241994016286177
$ws.Range("T15").Value = 0.002241994016286177

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 23.52144266666667
$ws.Range("H16").Value = 70.564328
$ws.Range("I16").Value = 0.01008349266740757
$ws.Range("J16").Value = 0.01008349266740757
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 12.17451333333334
$ws.Range("N16").Value = 36.52354
$ws.Range("O16").Value = 0.4763944950360188
$ws.Range("P16").Value = 0.4763944950360188
$ws.Range("Q16").Value = 286.362117364569
$ws.Range("R16").Value = 2577.259056281121
$ws.Range("S16").Value = 0.004803720397489029
$ws.Range("T16").Value = 0.004803720397489027

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 23.52144266666667
$ws.Range("H17").Value = 70.564328
$ws.Range("I17").Value = 0.01008349266740757
$ws.Range("J17").Value = 0.01008349266740757
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.300086333333334
$ws.Range("N17").Value = 15.900259
$ws.Range("O17").Value = 0.2073948981190463
$ws.Range("P17").Value = 0.2073948981190463
$ws.Range("Q17").Value = 124.6656768178836
$ws.Range("R17").Value = 1121.991091360952
$ws.Range("S17").Value = 0.002091264934441144
$ws.Range("T17").Value = 0.002091264934441143
